$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 153
$ws.Range("A153").Value = 2487632
$ws.Range("B153").Value = "C. D. MENSAJERO ISLA DE LA PALMA"
$ws.Range("C153").Value = "D. RODRIGUEZ GARCIA | J. RIES | M. NIANG | O. PEÑA LOPEZ | P. RODRIGUEZ RIVERO"
$ws.Range("D153").Value = 5
$ws.Range("E153").Value = 98
$ws.Range("F153").Value = 1.63
$ws.Range("G153").Value = 2
$ws.Range("H153").Value = 4
$ws.Range("I153").Value = 3.88
$ws.Range("J153").Value = 4
$ws.Range("K153").Value = 51.55
$ws.Range("L153").Value = 100
$ws.Range("M153").Value = -48.45
$ws.Range("N153").Value = 4
$ws.Range("O153").Value = 2
$ws.Range("P153").Value = 1
$ws.Range("Q153").Value = 2
$ws.Range("R153").Value = 5
$ws.Range("S153").Value = 0
$ws.Range("T153").Value = 2
$ws.Range("U153").Value = 3
$ws.Range("V153").Value = "Liga Regular `"B-B`""
$ws.Range("W153").Value = 8
$ws.Range("X153").Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"

# Row 154
$ws.Range("A154").Value = 2487632
$ws.Range("B154").Value = "C. D. MENSAJERO ISLA DE LA PALMA"
$ws.Range("C154").Value = "A. APARICIO IZQUIERDO | D. RODRIGUEZ GARCIA | J. RIES | O. PEÑA LOPEZ | P. RODRIGUEZ RIVERO"
$ws.Range("D154").Value = 5
$ws.Range("E154").Value = 27
$ws.Range("F154").Value = 0.45
$ws.Range("G154").Value = 2
$ws.Range("H154").Value = 3
$ws.Range("I154").Value = 0.88
$ws.Range("J154").Value = 1
$ws.Range("K154").Value = 227.27
$ws.Range("L154").Value = 300
$ws.Range("M154").Value = -72.73
$ws.Range("N154").Value = 0
$ws.Range("O154").Value = 2
$ws.Range("P154").Value = 0
$ws.Range("Q154").Value = 0
$ws.Range("R154").Value = 1
$ws.Range("S154").Value = 0
$ws.Range("T154").Value = 0
$ws.Range("U154").Value = 0
$ws.Range("V154").Value = "Liga Regular `"B-B`""
$ws.Range("W154").Value = 8
$ws.Range("X154").Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"

# Row 155
$ws.Range("A155").Value = 2487632
$ws.Range("B155").Value = "C.B. TRES CANTOS"
$ws.Range("C155").Value = "D. GONZALEZ LONGARELA | F. GOMEZ DE ENTERRIA LOPEZ | G. DIAZ MONTERO | J. DOMINGUEZ LARRE | N. MAIGA"
$ws.Range("D155").Value = 5
$ws.Range("E155").Value = 60
$ws.Range("F155").Value = 1
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 2
$ws.Range("I155").Value = 2
$ws.Range("J155").Value = 1
$ws.Range("K155").Value = 0
$ws.Range("L155").Value = 200
$ws.Range("M155").Value = -200
$ws.Range("N155").Value = 3
$ws.Range("O155").Value = 0
$ws.Range("P155").Value = 2
$ws.Range("Q155").Value = 3
$ws.Range("R155").Value = 1
$ws.Range("S155").Value = 0
$ws.Range("T155").Value = 0
$ws.Range("U155").Value = 0
$ws.Range("V155").Value = "Liga Regular `"B-B`""
$ws.Range("W155").Value = 8
$ws.Range("X155").Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"

# Row 156
$ws.Range("A156").Value = 2487632
$ws.Range("B156").Value = "C.B. TRES CANTOS"
$ws.Range("C156").Value = "D. GONZALEZ LONGARELA | G. DIAZ MONTERO | J. ATIENZA PEREA | J. DOMINGUEZ LARRE | N. MAIGA"
$ws.Range("D156").Value = 5
$ws.Range("E156").Value = 38
$ws.Range("F156").Value = 0.63
$ws.Range("G156").Value = 4
$ws.Range("H156").Value = 0
$ws.Range("I156").Value = 2
$ws.Range("J156").Value = 2.88
$ws.Range("K156").Value = 200
$ws.Range("L156").Value = 0
$ws.Range("M156").Value = 200
$ws.Range("N156").Value = 2
$ws.Range("O156").Value = 0
$ws.Range("P156").Value = 0
$ws.Range("Q156").Value = 0
$ws.Range("R156").Value = 3
$ws.Range("S156").Value = 2
$ws.Range("T156").Value = 1
$ws.Range("U156").Value = 2
$ws.Range("V156").Value = "Liga Regular `"B-B`""
$ws.Range("W156").Value = 8
$ws.Range("X156").Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"

# Row 157
$ws.Range("A157").Value = 2487632
$ws.Range("B157").Value = "C.B. TRES CANTOS"
$ws.Range("C157").Value = "A. SANCHO PEREZ | G. DIAZ MONTERO | J. ATIENZA PEREA | J. DOMINGUEZ LARRE | N. MAIGA"
$ws.Range("D157").Value = 5
$ws.Range("E157").Value = 27
$ws.Range("F157").Value = 0.45
$ws.Range("G157").Value = 3
$ws.Range("H157").Value = 2
$ws.Range("I157").Value = 1
$ws.Range("J157").Value = 0.88
$ws.Range("K157").Value = 300
$ws.Range("L157").Value = 227.27
$ws.Range("M157").Value = 72.73
$ws.Range("N157").Value = 1
$ws.Range("O157").Value = 0
$ws.Range("P157").Value = 0
$ws.Range("Q157").Value = 0
$ws.Range("R157").Value = 0
$ws.Range("S157").Value = 2
$ws.Range("T157").Value = 0
$ws.Range("U157").Value = 0
$ws.Range("V157").Value = "Liga Regular `"B-B`""
$ws.Range("W157").Value = 8
$ws.Range("X157").Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
